{"js": "// \"Version 2.\" -> \"Version 1.\"\n//\n// Scope the search/replace to just the lone \"2\" rather than the whole\n// \"Version 2.\" string: the digit sits immediately before a bookmark\n// (_GoBack) in this document, and replacing a range that *spans* the\n// bookmark would delete it. Matching only the digit (as a whole word, so\n// we don't clobber part of some other number) and swapping it for \"1\"\n// leaves every other run / the proofErr spell-check markers / the\n// bookmark untouched - exactly like someone retyping just that one\n// character would.\nconst body = context.document.body;\n\nconst results = body.search(\"2\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"1\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"Version 2.\" -> \"Version 1.\"\n#\n# Scope the Find/Replace to just the digit itself rather than the whole\n# \"Version 2.\" string. The digit sits immediately before a bookmark\n# (_GoBack) in this document, and replacing a range that *spans* the\n# bookmark would delete it - so, just like a person retyping only the\n# \"2\", we touch only that single character and leave everything else\n# (runs, proofErr spell-check markers, bookmark) alone.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"2\"\n$find.Replacement.Text = \"1\"\n\n$find.Execute(\n    [ref]\"2\",     # FindText\n    $false,       # MatchCase\n    $true,        # MatchWholeWord\n    $false,       # MatchWildcards\n    $false,       # MatchSoundsLike\n    $false,       # MatchAllWordForms\n    $true,        # Forward\n    1,            # Wrap (wdFindContinue)\n    $false,       # Format\n    [ref]\"1\",     # ReplaceWith\n    1             # Replace (wdReplaceOne)\n) | Out-Null\n"}
